# The underlying commit swaps the XML content of ppt/theme/theme1.xml
# (used by the slide master -> "Integral" / Red Violet) and
# ppt/theme/theme2.xml (used by the notes master -> default "Office Theme"),
# so that the deck's visible theme becomes the standard Office colour
# scheme.  The font scheme and format scheme are identical between the two
# theme parts - only the 12 colour-scheme slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) actually differ, so re-pointing the
# presentation's theme colours at the Office palette reproduces the
# substantive effect of the swap.
#
# ThemeColorScheme indices map 1:1 onto the OOXML <a:clrScheme> children in
# document order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
# 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink. COM RGB values are
# packed as r + g*256 + b*65536 (standard COLORREF/BGR ordering).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      -> #000000
$cs.Item(2).RGB  = 16777215   # lt1      -> #FFFFFF
$cs.Item(3).RGB  = 6968388     # dk2      -> #44546A
$cs.Item(4).RGB  = 15132391   # lt2      -> #E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  -> #5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  -> #ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  -> #A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  -> #FFC000
$cs.Item(9).RGB  = 12874308   # accent5  -> #4472C4
$cs.Item(10).RGB = 4697456    # accent6  -> #70AD47
$cs.Item(11).RGB = 12673797   # hlink    -> #0563C1
$cs.Item(12).RGB = 7491477    # folHlink -> #954F72
